$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "복소수 기초"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/01/05/complex_number_basic.html"

$ws.Range("D9").Value = "국내 대기업 & 정부기관 “딥러닝” 교육 실태? 작태?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/korean-deep-learning-quality/#utm_source=rss&utm_medium=rss&utm_campaign=korean-deep-learning-quality"

$ws.Range("D28").Value = "[임피던스 제어] Improving Low-Impedance Performance (Series Dynamics)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/173"

$ws.Range("D46").Value = "COVID-19 백신 부작용의 시간별 특징"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/423"
